# The document body consists of a single, empty paragraph. The edit turns
# that empty paragraph into one containing a single run with the text "12".
$d = $word.ActiveDocument

$p = $d.Paragraphs.Item(1)
$p.Range.Text = "12"
